$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Intro paragraph: add "I" before "developed", and "and implement" before "an A*" ---
Replace-Text "this project, developed a Java program" "this project, I developed a Java program"
Replace-Text "we were required to develop an A* algorithm" "we were required to develop and implement an A* algorithm"

# --- Insert a blank paragraph before "The A* Algorithm" heading ---
$rng = $d.Content
$rng.Find.Execute("The A* Algorithm") | Out-Null
$headingRange = $rng.Duplicate
$headingRange.Collapse(1)
$headingRange.InsertParagraphBefore()

# --- COST section ---
Replace-Text "This was equal to the number of steps to reach the goal state." "This is equal to the number of steps required to reach the goal state."

# --- UNIT TESTING section ---
Replace-Text "JUnit Test were implemented in Eclipse for several test cases." "JUnit Testing is implemented in Eclipse for several test cases."
Replace-Text "The purpose of the tests was to determine if suboptimal solutions were being found." "The purpose of the tests is to determine if suboptimal solutions exist."
Replace-Text "For example, for the case of Pitchers={2,5,6,72}, goal={143}, the algorithm returns 8, when the correct solution is 7." "For example, is Pitchers={2,5,6,72} and goal={143}, the algorithm returns 8, when the optimal solution is 7."

# --- Remove the trailing empty ListParagraph paragraph at the end of the document ---
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$tailRange = $d.Range($secondLast.Range.End - 1, $d.Content.End)
$tailRange.Delete()
